{"js": "// Replace the 25 multiplication-problem texts in the single table on the\n// page with their new values, preserving all run/paragraph formatting\n// (font, size, justification) by rewriting only the text inside each\n// cell's existing paragraph range.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col) -> new text, taken from the canonical OOXML diff.\nconst updates = [\n  [0, 0, \"17\u00d765=\"],\n  [0, 1, \"28\u00d785=\"],\n  [0, 2, \"84\u00d755=\"],\n  [0, 3, \"39\u00d774=\"],\n  [0, 4, \"47\u00d756=\"],\n\n  [4, 0, \"37\u00d762=\"],\n  [4, 1, \"15\u00d792=\"],\n  [4, 2, \"22\u00d792=\"],\n  [4, 3, \"16\u00d758=\"],\n  [4, 4, \"31\u00d784=\"],\n\n  [9, 0, \"72\u00d721=\"],\n  [9, 1, \"16\u00d745=\"],\n  [9, 2, \"29\u00d797=\"],\n  [9, 3, \"41\u00d737=\"],\n  [9, 4, \"52\u00d718=\"],\n\n  [14, 0, \"66\u00d758=\"],\n  [14, 1, \"39\u00d747=\"],\n  [14, 2, \"32\u00d795=\"],\n  [14, 3, \"20\u00d769=\"],\n  [14, 4, \"33\u00d718=\"],\n\n  [19, 0, \"54\u00d756=\"],\n  [19, 1, \"25\u00d730=\"],\n  [19, 2, \"34\u00d747=\"],\n  [19, 3, \"27\u00d726=\"],\n  [19, 4, \"88\u00d792=\"],\n];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(text, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 multiplication-problem texts in the single table on the\n# page with their new values, preserving all run/paragraph formatting\n# (font, size, justification) by overwriting the text of each cell's\n# Range directly (Word re-uses the existing run formatting when the\n# Range.Text is simply reassigned).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based (row, col) -> new text, taken from the canonical OOXML diff.\n$updates = @(\n    @(1, 1, \"17\u00d765=\"),\n    @(1, 2, \"28\u00d785=\"),\n    @(1, 3, \"84\u00d755=\"),\n    @(1, 4, \"39\u00d774=\"),\n    @(1, 5, \"47\u00d756=\"),\n\n    @(5, 1, \"37\u00d762=\"),\n    @(5, 2, \"15\u00d792=\"),\n    @(5, 3, \"22\u00d792=\"),\n    @(5, 4, \"16\u00d758=\"),\n    @(5, 5, \"31\u00d784=\"),\n\n    @(10, 1, \"72\u00d721=\"),\n    @(10, 2, \"16\u00d745=\"),\n    @(10, 3, \"29\u00d797=\"),\n    @(10, 4, \"41\u00d737=\"),\n    @(10, 5, \"52\u00d718=\"),\n\n    @(15, 1, \"66\u00d758=\"),\n    @(15, 2, \"39\u00d747=\"),\n    @(15, 3, \"32\u00d795=\"),\n    @(15, 4, \"20\u00d769=\"),\n    @(15, 5, \"33\u00d718=\"),\n\n    @(20, 1, \"54\u00d756=\"),\n    @(20, 2, \"25\u00d730=\"),\n    @(20, 3, \"34\u00d747=\"),\n    @(20, 4, \"27\u00d726=\"),\n    @(20, 5, \"88\u00d792=\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $text\n}\n"}
